$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look like plain numbers as Text,
# so Excel keeps them as strings (matching the source inlineStr cells)
# instead of silently coercing them into floating-point numbers.
$numericLookingCells = @('D5','D6','D8','D10','D12','D13','D14','D19','D20','D21','D22','D25','D26','D27','D28','D31','D32','D34','D36','D37','D38','D39','D40','D41','D43','D44','D46','D47','D48','D49','D50','D51')
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.695.08'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '2.406.78'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '564.93'
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('D6').Value = '137.85'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = '0.537'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').Value = '2.388.98'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('D10').Value = '0.105'
$ws.Range('E10').Value = '  -3.07%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '5.03'
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('D13').Value = '0.334'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').Value = '25.71'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '2.848.16'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = '60.798.06'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = '2.390.37'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('D19').Value = '8.00'
$ws.Range('E19').Value = '  +11.03%  '
$ws.Range('D20').Value = '10.43'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').Value = '321.59'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').Value = '4.02'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '1.80'
$ws.Range('E25').Value = '  -5.37%  '
$ws.Range('D26').Value = '64.19'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').Value = '564.68'
$ws.Range('E27').Value = '  -3.20%  '
$ws.Range('D28').Value = '8.18'
$ws.Range('E28').Value = '  -10.95%  '
$ws.Range('D29').Value = '2.531.83'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '0.0₃0914'
$ws.Range('E30').Value = '  -2.07%  '
$ws.Range('D31').Value = '7.89'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '1.31'
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('E33').Value = '  -3.77%  '
$ws.Range('D34').Value = '0.131'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '152.35'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.39'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Value = '0.366'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '4.50'
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('D40').Value = '18.10'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('D41').Value = '5.06'
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '1.64'
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('D44').Value = '2.34'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').Value = '0.0₆0286'
$ws.Range('E45').Value = '  +3.33%  '
$ws.Range('D46').Value = '141.48'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').Value = '3.49'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('D48').Value = '0.582'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '0.0496'
$ws.Range('E49').Value = '  -2.60%  '
$ws.Range('D50').Value = '18.98'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').Value = '0.0898'
$ws.Range('E51').Value = '  +0.18%  '
